$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:71 down to 60:72
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 10
$ws.Cells.Item(59, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(59, 3).Value = "La Araucanía"
$ws.Cells.Item(59, 4).Value = 44985
$ws.Cells.Item(59, 5).Value = 9
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100108
$ws.Cells.Item(59, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(59, 9).Value = 100108003
$ws.Cells.Item(59, 10).Value = "Maracuyá"
$ws.Cells.Item(59, 11).Value = "Sin especificar"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 15
$ws.Cells.Item(59, 14).Value = 65000
$ws.Cells.Item(59, 15).Value = 65000
$ws.Cells.Item(59, 16).Value = 65000
$ws.Cells.Item(59, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(59, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 19).Value = 3611
$ws.Cells.Item(59, 20).Value = 18

# Ensure date-formatted column D uses the same number format as the rest of the column
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
